$d = $word.ActiveDocument

$d.Content.Find.Execute("24×86=2064", $true, $false, $false, $false, $false, $true, 1, $false, "65×15=975", 2) | Out-Null
$d.Content.Find.Execute("20×21=420", $true, $false, $false, $false, $false, $true, 1, $false, "36×12=432", 2) | Out-Null
$d.Content.Find.Execute("85×65=5525", $true, $false, $false, $false, $false, $true, 1, $false, "37×17=629", 2) | Out-Null
$d.Content.Find.Execute("51×97=4947", $true, $false, $false, $false, $false, $true, 1, $false, "15×93=1395", 2) | Out-Null
$d.Content.Find.Execute("30×78=2340", $true, $false, $false, $false, $false, $true, 1, $false, "51×81=4131", 2) | Out-Null
$d.Content.Find.Execute("98×85=8330", $true, $false, $false, $false, $false, $true, 1, $false, "16×40=640", 2) | Out-Null
$d.Content.Find.Execute("86×67=5762", $true, $false, $false, $false, $false, $true, 1, $false, "19×33=627", 2) | Out-Null
$d.Content.Find.Execute("18×15=270", $true, $false, $false, $false, $false, $true, 1, $false, "56×65=3640", 2) | Out-Null
$d.Content.Find.Execute("99×76=7524", $true, $false, $false, $false, $false, $true, 1, $false, "35×98=3430", 2) | Out-Null
$d.Content.Find.Execute("91×14=1274", $true, $false, $false, $false, $false, $true, 1, $false, "71×44=3124", 2) | Out-Null
$d.Content.Find.Execute("54×57=3078", $true, $false, $false, $false, $false, $true, 1, $false, "65×63=4095", 2) | Out-Null
$d.Content.Find.Execute("59×62=3658", $true, $false, $false, $false, $false, $true, 1, $false, "93×19=1767", 2) | Out-Null
$d.Content.Find.Execute("39×12=468", $true, $false, $false, $false, $false, $true, 1, $false, "90×79=7110", 2) | Out-Null
$d.Content.Find.Execute("72×91=6552", $true, $false, $false, $false, $false, $true, 1, $false, "95×21=1995", 2) | Out-Null
$d.Content.Find.Execute("43×24=1032", $true, $false, $false, $false, $false, $true, 1, $false, "94×43=4042", 2) | Out-Null
$d.Content.Find.Execute("71×47=3337", $true, $false, $false, $false, $false, $true, 1, $false, "56×39=2184", 2) | Out-Null
$d.Content.Find.Execute("57×49=2793", $true, $false, $false, $false, $false, $true, 1, $false, "59×66=3894", 2) | Out-Null
$d.Content.Find.Execute("85×45=3825", $true, $false, $false, $false, $false, $true, 1, $false, "60×54=3240", 2) | Out-Null
$d.Content.Find.Execute("44×69=3036", $true, $false, $false, $false, $false, $true, 1, $false, "34×99=3366", 2) | Out-Null
$d.Content.Find.Execute("98×27=2646", $true, $false, $false, $false, $false, $true, 1, $false, "61×33=2013", 2) | Out-Null
$d.Content.Find.Execute("52×60=3120", $true, $false, $false, $false, $false, $true, 1, $false, "34×57=1938", 2) | Out-Null
$d.Content.Find.Execute("81×65=5265", $true, $false, $false, $false, $false, $true, 1, $false, "15×14=210", 2) | Out-Null
$d.Content.Find.Execute("32×81=2592", $true, $false, $false, $false, $false, $true, 1, $false, "13×26=338", 2) | Out-Null
$d.Content.Find.Execute("88×39=3432", $true, $false, $false, $false, $false, $true, 1, $false, "46×98=4508", 2) | Out-Null
$d.Content.Find.Execute("35×61=2135", $true, $false, $false, $false, $false, $true, 1, $false, "72×34=2448", 2) | Out-Null
